$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.511.97"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.344.50"
$ws.Range("E3").Value = "  -4.30%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "535.87"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").Value = "135.81"
$ws.Range("E6").Value = "  -6.57%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  -11.47%  "
$ws.Range("D9").Value = "2.340.82"
$ws.Range("E9").Value = "  -4.49%  "
$ws.Range("D10").Value = "0.103"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "5.21"
$ws.Range("E12").Value = "  -3.49%  "
$ws.Range("D13").Value = "0.338"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").Value = "24.48"
$ws.Range("E14").Value = "  -5.31%  "
$ws.Range("D15").Value = "2.766.34"
$ws.Range("E15").Value = "  -4.31%  "
$ws.Range("D16").Value = "60.418.26"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "0.0000160"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "2.348.53"
$ws.Range("E18").Value = "  -4.32%  "
$ws.Range("D19").Value = "10.48"
$ws.Range("E19").Value = "  -4.63%  "
$ws.Range("D20").Value = "312.79"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").Value = "4.02"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "6.50"
$ws.Range("E22").Value = "  -6.35%  "
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").Value = "1.86"
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "62.72"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "8.32"
$ws.Range("E27").Value = "  +10.53%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.460.68"
$ws.Range("E29").Value = "  -4.63%  "
$ws.Range("D30").Value = "7.84"
$ws.Range("E30").Value = "  -3.45%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0879"
$ws.Range("E31").Value = "  -7.81%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "495.99"
$ws.Range("E32").Value = "  -6.90%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  -5.45%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.143"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "1.76"
$ws.Range("E35").Value = "  -6.39%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.51"
$ws.Range("E36").Value = "  -3.58%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.50"
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("D39").Value = "0.369"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.18"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "5.17"
$ws.Range("E41").Value = "  -9.97%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.75"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "138.81"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "39.88"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  -7.90%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "136.63"
$ws.Range("E47").Value = "  -4.82%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "3.47"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0504"
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "19.26"
$ws.Range("E50").Value = "  -9.10%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.563"
$ws.Range("E51").Value = "  -3.68%  "
